$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$ws.Range("A5").Value = "eng.caio.camilo@gmail.com"
$ws.Range("A5").Hyperlink = "mailto:eng.caio.camilo@gmail.com"
$ws.Range("A5").Select()
